$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "reviews_count" column (E) is being removed from the data set; delete
# the entire column so every later column (reviews_average, latitude,
# longitude, is_permanently_closed, gmaps_link, latest_review_date) shifts
# one place to the left (F->E, G->F, H->G, I->H, J->I, K->J).
$ws.Range("E1").EntireColumn.Delete()
